$d = $word.ActiveDocument
try {
  $d.CopyStylesFromTemplate("C:\Users\admin\Documents\test1.docx")
  Write-Output "ok"
} catch {
  Write-Output "ERROR: $_"
}
